# Fixed minor BOM errors
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 21: Name was "Thermistor Connector" -> should read the JST connector's own part number
$ws.Range("A21").Value = "B2B-EH-A (LF)(SN)"

# Row 9 (CL21A475KAQNNNG 4.7uF cap): quantity correction 3 -> 2, subtotal recalculated
$ws.Range("D9").Value = 2
$ws.Range("H9").Value = 0.192

# Row 11 (GRM21BR61E106KA73K 10uF cap): quantity correction 2 -> 5, subtotal recalculated
$ws.Range("D11").Value = 5
$ws.Range("H11").Value = 0.8175

# Row 13 (0402B682K500CT 680pF cap): subtotal price fix
$ws.Range("H13").Value = 0.01

# Row 29 (FTSH-105-01-L-DV-K): unit price / subtotal price fix
$ws.Range("G29").Value = 2.85
$ws.Range("H29").Value = 2.85

# Row 49 (TPS55340 boost converter): correct part number typo and pricing
$ws.Range("A49").Value = "TPS55340QRTERQ1"
$ws.Range("C49").Value = "TPS55340QRTERQ1"
$ws.Range("G49").Value = 3.96
$ws.Range("H49").Value = 3.96

$wb.Save()
